# Add a new "Greece" market sheet, cloned from the existing "Croatia" sheet,
# matching the layout/styles/merges/column widths of the other market sheets.

$wb = $excel.ActiveWorkbook

# Clear the current selection state on "Croatia" first (full-sheet select,
# i.e. what's left behind on the template sheet once it stops being the
# active tab) before we clone it, so the clone inherits the old B7
# selection and Croatia ends up with the "select all" state.
$croatia = $wb.Worksheets.Item("Croatia")
$croatia.Cells.Select() | Out-Null

# Clone "Croatia" (preserves column widths, styles, merged cells, page
# setup, etc.) and drop the copy right after it.
$croatia.Copy($null, $croatia)
$greece = $wb.Worksheets.Item("Croatia (2)")
$greece.Name = "Greece"

# Greece-specific data: part number (B4) then market name (B2), matching
# the authoring order of the new shared-string entries.
$greece.Range("B4").Value = "NGC-4119/T3190"
$greece.Range("B2").Value = "Greece Market"

# Leave the new sheet as the active tab/selection.
$greece.Range("C13").Select() | Out-Null
